$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.220.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "'1.797.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.57%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'314.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("D7").Value = "'0.5190"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.58%  "

$ws.Range("D8").Value = "'0.3818"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.28%  "

$ws.Range("D9").Value = "'0.07923"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.33%  "

$ws.Range("D10").Value = "'41.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("E11").Value = "  -1.35%  "

$ws.Range("D12").Value = "'6.277"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.35%  "

$ws.Range("D13").Value = "'1.002"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("E14").Value = "  -2.73%  "

$ws.Range("D15").Value = "'1.792.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.53%  "

$ws.Range("D16").Value = "'7.268"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.94%  "

$ws.Range("D17").Value = "'93.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("E18").Value = "  -3.26%  "

$ws.Range("D19").Value = "'0.06560"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.52%  "

$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("E21").Value = "  -2.93%  "

$ws.Range("D22").Value = "'5.954"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.58%  "

$ws.Range("D23").Value = "'28.243.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.16%  "

$ws.Range("D24").Value = "'11.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.36%  "

$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("D26").Value = "'160.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.32%  "

$ws.Range("D27").Value = "'20.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.05%  "

$ws.Range("D28").Value = "'2.000.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.45%  "

$ws.Range("D29").Value = "'2.338"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.12%  "

$ws.Range("D30").Value = "'123.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.65%  "

$ws.Range("D31").Value = "'0.1067"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.89%  "

$ws.Range("D32").Value = "'1.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.49%  "

$ws.Range("D33").Value = "'3.671"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("D34").Value = "'5.568"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.50%  "

$ws.Range("D35").Value = "'0.07323"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.63%  "

$ws.Range("D36").Value = "'12.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.59%  "

$ws.Range("D37").Value = "'0.02331"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.08%  "

$ws.Range("D38").Value = "'0.2142"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.05%  "

$ws.Range("D39").Value = "'5.070"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.52%  "

$ws.Range("D40").Value = "'8.617"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("D41").Value = "'0.6168"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.09%  "

$ws.Range("E42").Value = "  -1.68%  "

$ws.Range("D43").Value = "'1.371"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.05%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.6041"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.11%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.11%  "

$ws.Range("D46").Value = "'3.785"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.17%  "

$ws.Range("D47").Value = "'127.49"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").Value = "'1.232"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").Value = "'1.924"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.50%  "

$ws.Range("D50").Value = "'0.06782"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.36%  "

$ws.Range("D51").Value = "'73.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.60%  "
